$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,17

$data[0,0] = 'ECs'
$data[0,1] = 3
$data[0,2] = 1
$data[0,3] = 9.031965666666666
$data[0,4] = 27.095897
$data[0,5] = 0.4424406034784756
$data[0,6] = 0.4424406034784755
$data[0,7] = 3
$data[0,8] = 1
$data[0,9] = 8.540560666666666
$data[0,10] = 25.621682
$data[0,11] = 0.4159358086620884
$data[0,12] = 0.4159358086620884
$data[0,13] = 77.13805071541711
$data[0,14] = 694.2424564387539
$data[0,15] = 0.1840268901927621
$data[0,16] = 0.1840268901927621

$data[1,0] = 'FAPs'
$data[1,1] = 3
$data[1,2] = 1
$data[1,3] = 9.031965666666666
$data[1,4] = 27.095897
$data[1,5] = 0.4424406034784756
$data[1,6] = 0.4424406034784755
$data[1,7] = 3
$data[1,8] = 1
$data[1,9] = 11.57455166666667
$data[1,10] = 34.723655
$data[1,11] = 0.563694901924408
$data[1,12] = 0.563694901924408
$data[1,13] = 104.5409532603928
$data[1,14] = 940.868579343535
$data[1,15] = 0.2494015125851752
$data[1,16] = 0.2494015125851751

$data[2,0] = 'Inflammatory-Mac'
$data[2,1] = 3
$data[2,2] = 1
$data[2,3] = 9.031965666666666
$data[2,4] = 27.095897
$data[2,5] = 0.4424406034784756
$data[2,6] = 0.4424406034784755
$data[2,7] = 1
$data[2,8] = 0.3333333333333333
$data[2,9] = 0.026642
$data[2,10] = 0.079926
$data[2,11] = 0.001297498167494471
$data[2,12] = 0.001297498167494471
$data[2,13] = 0.2406296292913333
$data[2,14] = 2.165666663622
$data[2,15] = 0.00057406587223847
$data[2,16] = 0.0005740658722384699

$data[3,0] = 'MuSCs'
$data[3,1] = 3
$data[3,2] = 1
$data[3,3] = 9.031965666666666
$data[3,4] = 27.095897
$data[3,5] = 0.4424406034784756
$data[3,6] = 0.4424406034784755
$data[3,7] = 3
$data[3,8] = 1
$data[3,9] = 0.391608
$data[3,10] = 1.174824
$data[3,11] = 0.01907179124600912
$data[3,12] = 0.01907179124600912
$data[3,13] = 3.536990010792
$data[3,14] = 31.832910097128
$data[3,15] = 0.008438134828299781
$data[3,16] = 0.008438134828299781

$data[4,0] = 'ECs'
$data[4,1] = 3
$data[4,2] = 1
$data[4,3] = 9.124904999999998
$data[4,4] = 27.37471499999999
$data[4,5] = 0.4469933372071527
$data[4,6] = 0.4469933372071526
$data[4,7] = 3
$data[4,8] = 1
$data[4,9] = 8.540560666666666
$data[4,10] = 25.621682
$data[4,11] = 0.4159358086620884
$data[4,12] = 0.4159358086620884
$data[4,13] = 77.93180473006998
$data[4,14] = 701.3862425706299
$data[4,15] = 0.1859205351778226
$data[4,16] = 0.1859205351778226

$data[5,0] = 'FAPs'
$data[5,1] = 3
$data[5,2] = 1
$data[5,3] = 9.124904999999998
$data[5,4] = 27.37471499999999
$data[5,5] = 0.4469933372071527
$data[5,6] = 0.4469933372071526
$data[5,7] = 3
$data[5,8] = 1
$data[5,9] = 11.57455166666667
$data[5,10] = 34.723655
$data[5,11] = 0.563694901924408
$data[5,12] = 0.563694901924408
$data[5,13] = 105.616684375925
$data[5,14] = 950.5501593833249
$data[5,15] = 0.2519678653778498
$data[5,16] = 0.2519678653778497

$data[6,0] = 'Inflammatory-Mac'
$data[6,1] = 3
$data[6,2] = 1
$data[6,3] = 9.124904999999998
$data[6,4] = 27.37471499999999
$data[6,5] = 0.4469933372071527
$data[6,6] = 0.4469933372071526
$data[6,7] = 1
$data[6,8] = 0.3333333333333333
$data[6,9] = 0.026642
$data[6,10] = 0.079926
$data[6,11] = 0.001297498167494471
$data[6,12] = 0.001297498167494471
$data[6,13] = 0.2431057190099999
$data[6,14] = 2.187951471089999
$data[6,15] = 0.0005799730359085188
$data[6,16] = 0.0005799730359085188

$data[7,0] = 'MuSCs'
$data[7,1] = 3
$data[7,2] = 1
$data[7,3] = 9.124904999999998
$data[7,4] = 27.37471499999999
$data[7,5] = 0.4469933372071527
$data[7,6] = 0.4469933372071526
$data[7,7] = 3
$data[7,8] = 1
$data[7,9] = 0.391608
$data[7,10] = 1.174824
$data[7,11] = 0.01907179124600912
$data[7,12] = 0.01907179124600912
$data[7,13] = 3.573385797239999
$data[7,14] = 32.16047217516
$data[7,15] = 0.008524963615571776
$data[7,16] = 0.008524963615571774

$data[8,0] = 'ECs'
$data[8,1] = 2
$data[8,2] = 0.6666666666666666
$data[8,3] = 0.2635683333333333
$data[8,4] = 0.790705
$data[8,5] = 0.01291117977653399
$data[8,6] = 0.01291117977653399
$data[8,7] = 3
$data[8,8] = 1
$data[8,9] = 8.540560666666666
$data[8,10] = 25.621682
$data[8,11] = 0.4159358086620884
$data[8,12] = 0.4159358086620884
$data[8,13] = 2.251021340645555
$data[8,14] = 20.25919206581
$data[8,15] = 0.005370222001134268
$data[8,16] = 0.005370222001134267

$data[9,0] = 'FAPs'
$data[9,1] = 2
$data[9,2] = 0.6666666666666666
$data[9,3] = 0.2635683333333333
$data[9,4] = 0.790705
$data[9,5] = 0.01291117977653399
$data[9,6] = 0.01291117977653399
$data[9,7] = 3
$data[9,8] = 1
$data[9,9] = 11.57455166666667
$data[9,10] = 34.723655
$data[9,11] = 0.563694901924408
$data[9,12] = 0.563694901924408
$data[9,13] = 3.050685291863889
$data[9,14] = 27.456167626775
$data[9,15] = 0.007277966217861729
$data[9,16] = 0.007277966217861727

$data[10,0] = 'Inflammatory-Mac'
$data[10,1] = 2
$data[10,2] = 0.6666666666666666
$data[10,3] = 0.2635683333333333
$data[10,4] = 0.790705
$data[10,5] = 0.01291117977653399
$data[10,6] = 0.01291117977653399
$data[10,7] = 1
$data[10,8] = 0.3333333333333333
$data[10,9] = 0.026642
$data[10,10] = 0.079926
$data[10,11] = 0.001297498167494471
$data[10,12] = 0.001297498167494471
$data[10,13] = 0.007021987536666667
$data[10,14] = 0.06319788782999999
$data[10,15] = 0.00001675223210024453
$data[10,16] = 0.00001675223210024453

$data[11,0] = 'MuSCs'
$data[11,1] = 2
$data[11,2] = 0.6666666666666666
$data[11,3] = 0.2635683333333333
$data[11,4] = 0.790705
$data[11,5] = 0.01291117977653399
$data[11,6] = 0.01291117977653399
$data[11,7] = 3
$data[11,8] = 1
$data[11,9] = 0.391608
$data[11,10] = 1.174824
$data[11,11] = 0.01907179124600912
$data[11,12] = 0.01907179124600912
$data[11,13] = 0.10321546788
$data[11,14] = 0.9289392109200001
$data[11,15] = 0.000246239325437751
$data[11,16] = 0.0002462393254377509

$data[12,0] = 'ECs'
$data[12,1] = 3
$data[12,2] = 1
$data[12,3] = 1.809602666666667
$data[12,4] = 5.428808
$data[12,5] = 0.08864534315615299
$data[12,6] = 0.08864534315615297
$data[12,7] = 3
$data[12,8] = 1
$data[12,9] = 8.540560666666666
$data[12,10] = 25.621682
$data[12,11] = 0.4159358086620884
$data[12,12] = 0.4159358086620884
$data[12,13] = 15.45502135722844
$data[12,14] = 139.095192215056
$data[12,15] = 0.03687077248978282
$data[12,16] = 0.03687077248978281

$data[13,0] = 'FAPs'
$data[13,1] = 3
$data[13,2] = 1
$data[13,3] = 1.809602666666667
$data[13,4] = 5.428808
$data[13,5] = 0.08864534315615299
$data[13,6] = 0.08864534315615297
$data[13,7] = 3
$data[13,8] = 1
$data[13,9] = 11.57455166666667
$data[13,10] = 34.723655
$data[13,11] = 0.563694901924408
$data[13,12] = 0.563694901924408
$data[13,13] = 20.94533956147111
$data[13,14] = 188.50805605324
$data[13,15] = 0.04996892801646315
$data[13,16] = 0.04996892801646313

$data[14,0] = 'Inflammatory-Mac'
$data[14,1] = 3
$data[14,2] = 1
$data[14,3] = 1.809602666666667
$data[14,4] = 5.428808
$data[14,5] = 0.08864534315615299
$data[14,6] = 0.08864534315615297
$data[14,7] = 1
$data[14,8] = 0.3333333333333333
$data[14,9] = 0.026642
$data[14,10] = 0.079926
$data[14,11] = 0.001297498167494471
$data[14,12] = 0.001297498167494471
$data[14,13] = 0.04821143424533333
$data[14,14] = 0.433902908208
$data[14,15] = 0.0001150171703020271
$data[14,16] = 0.000115017170302027

$data[15,0] = 'MuSCs'
$data[15,1] = 3
$data[15,2] = 1
$data[15,3] = 1.809602666666667
$data[15,4] = 5.428808
$data[15,5] = 0.08864534315615299
$data[15,6] = 0.08864534315615297
$data[15,7] = 3
$data[15,8] = 1
$data[15,9] = 0.391608
$data[15,10] = 1.174824
$data[15,11] = 0.01907179124600912
$data[15,12] = 0.01907179124600912
$data[15,13] = 0.708654881088
$data[15,14] = 6.377893929792001
$data[15,15] = 0.001690625479604993
$data[15,16] = 0.001690625479604992

$data[16,0] = 'ECs'
$data[16,1] = 2
$data[16,2] = 0.6666666666666666
$data[16,3] = 0.1839203333333334
$data[16,4] = 0.5517610000000001
$data[16,5] = 0.009009536381684918
$data[16,6] = 0.009009536381684917
$data[16,7] = 3
$data[16,8] = 1
$data[16,9] = 8.540560666666666
$data[16,10] = 25.621682
$data[16,11] = 0.4159358086620884
$data[16,12] = 0.4159358086620884
$data[16,13] = 1.570782764666889
$data[16,14] = 14.137044882002
$data[16,15] = 0.003747388800586623
$data[16,16] = 0.003747388800586622

$data[17,0] = 'FAPs'
$data[17,1] = 2
$data[17,2] = 0.6666666666666666
$data[17,3] = 0.1839203333333334
$data[17,4] = 0.5517610000000001
$data[17,5] = 0.009009536381684918
$data[17,6] = 0.009009536381684917
$data[17,7] = 3
$data[17,8] = 1
$data[17,9] = 11.57455166666667
$data[17,10] = 34.723655
$data[17,11] = 0.563694901924408
$data[17,12] = 0.563694901924408
$data[17,13] = 2.128795400717222
$data[17,14] = 19.159158606455
$data[17,15] = 0.005078629727058265
$data[17,16] = 0.005078629727058265

$data[18,0] = 'Inflammatory-Mac'
$data[18,1] = 2
$data[18,2] = 0.6666666666666666
$data[18,3] = 0.1839203333333334
$data[18,4] = 0.5517610000000001
$data[18,5] = 0.009009536381684918
$data[18,6] = 0.009009536381684917
$data[18,7] = 1
$data[18,8] = 0.3333333333333333
$data[18,9] = 0.026642
$data[18,10] = 0.079926
$data[18,11] = 0.001297498167494471
$data[18,12] = 0.001297498167494471
$data[18,13] = 0.004900005520666667
$data[18,14] = 0.044100049686
$data[18,15] = 0.00001168985694521095
$data[18,16] = 0.00001168985694521095

$data[19,0] = 'MuSCs'
$data[19,1] = 2
$data[19,2] = 0.6666666666666666
$data[19,3] = 0.1839203333333334
$data[19,4] = 0.5517610000000001
$data[19,5] = 0.009009536381684918
$data[19,6] = 0.009009536381684917
$data[19,7] = 3
$data[19,8] = 1
$data[19,9] = 0.391608
$data[19,10] = 1.174824
$data[19,11] = 0.01907179124600912
$data[19,12] = 0.01907179124600912
$data[19,13] = 0.07202467389600001
$data[19,14] = 0.6482220650640002
$data[19,15] = 0.0001718279970948191
$data[19,16] = 0.0001718279970948191

$ws.Range("D2:T21").Value = $data
Write-Output "applied"